# Rebuild the "plot_infos_pcr_snv" results table: reorder the metric
# columns (R², RMSE, Offset, Slope) and add a "Validação" (validation)
# row for each attribute (SST, PH, AT, FIRMEZA (N), UBS (%)), plus the
# updated/recomputed metric values for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Attribute'
$ws.Cells.Item(1, 2).Value = 'Y'
$ws.Cells.Item(1, 3).Value = 'R²'
$ws.Cells.Item(1, 4).Value = 'RMSE'
$ws.Cells.Item(1, 5).Value = 'Offset'
$ws.Cells.Item(1, 6).Value = 'Slope'

$ws.Cells.Item(2, 1).Value = 'SST'
$ws.Cells.Item(2, 2).Value = 'Referência'
$ws.Cells.Item(2, 3).Value = 0.725752312300771
$ws.Cells.Item(2, 4).Value = 1.363035532866461
$ws.Cells.Item(2, 5).Value = 3.817057673880066
$ws.Cells.Item(2, 6).Value = 0.7257523123007712

$ws.Cells.Item(3, 1).Value = 'SST'
$ws.Cells.Item(3, 2).Value = 'Predição'
$ws.Cells.Item(3, 3).Value = 0.6844083075715417
$ws.Cells.Item(3, 4).Value = 1.462172028352482
$ws.Cells.Item(3, 5).Value = 3.975964250659266
$ws.Cells.Item(3, 6).Value = 0.714341310023144

$ws.Cells.Item(4, 1).Value = 'SST'
$ws.Cells.Item(4, 2).Value = 'Validação'
$ws.Cells.Item(4, 3).Value = 0.6750277633089374
$ws.Cells.Item(4, 4).Value = 1.28067716930816
$ws.Cells.Item(4, 5).Value = 0.8794082431567105
$ws.Cells.Item(4, 6).Value = 0.9410054317098512

$ws.Cells.Item(5, 1).Value = 'PH'
$ws.Cells.Item(5, 2).Value = 'Referência'
$ws.Cells.Item(5, 3).Value = 0.2208274643071334
$ws.Cells.Item(5, 4).Value = 0.2764102973443819
$ws.Cells.Item(5, 5).Value = 2.572306780362552
$ws.Cells.Item(5, 6).Value = 0.2208274643071334

$ws.Cells.Item(6, 1).Value = 'PH'
$ws.Cells.Item(6, 2).Value = 'Predição'
$ws.Cells.Item(6, 3).Value = 0.1754672905536661
$ws.Cells.Item(6, 4).Value = 0.2843422165241026
$ws.Cells.Item(6, 5).Value = 2.640215549496488
$ws.Cells.Item(6, 6).Value = 0.2003502245197424

$ws.Cells.Item(7, 1).Value = 'PH'
$ws.Cells.Item(7, 2).Value = 'Validação'
$ws.Cells.Item(7, 3).Value = 0.3472123147864786
$ws.Cells.Item(7, 4).Value = 0.199461780340821
$ws.Cells.Item(7, 5).Value = -1.242848334414727
$ws.Cells.Item(7, 6).Value = 1.370099221676297

$ws.Cells.Item(8, 1).Value = 'AT'
$ws.Cells.Item(8, 2).Value = 'Referência'
$ws.Cells.Item(8, 3).Value = 0.4463262469582551
$ws.Cells.Item(8, 4).Value = 0.4160801037858976
$ws.Cells.Item(8, 5).Value = 0.6231076055231972
$ws.Cells.Item(8, 6).Value = 0.446326246958255

$ws.Cells.Item(9, 1).Value = 'AT'
$ws.Cells.Item(9, 2).Value = 'Predição'
$ws.Cells.Item(9, 3).Value = 0.3938521138578016
$ws.Cells.Item(9, 4).Value = 0.435350731439894
$ws.Cells.Item(9, 5).Value = 0.6482324452604232
$ws.Cells.Item(9, 6).Value = 0.4241068854878716

$ws.Cells.Item(10, 1).Value = 'AT'
$ws.Cells.Item(10, 2).Value = 'Validação'
$ws.Cells.Item(10, 3).Value = 0.611451678176929
$ws.Cells.Item(10, 4).Value = 0.3201363566043914
$ws.Cells.Item(10, 5).Value = -0.2354241571362763
$ws.Cells.Item(10, 6).Value = 1.169788316409431

$ws.Cells.Item(11, 1).Value = 'FIRMEZA (N)'
$ws.Cells.Item(11, 2).Value = 'Referência'
$ws.Cells.Item(11, 3).Value = 0.4351124807062247
$ws.Cells.Item(11, 4).Value = 75.97813136126791
$ws.Cells.Item(11, 5).Value = 290.425730044084
$ws.Cells.Item(11, 6).Value = 0.4351124807062243

$ws.Cells.Item(12, 1).Value = 'FIRMEZA (N)'
$ws.Cells.Item(12, 2).Value = 'Predição'
$ws.Cells.Item(12, 3).Value = 0.3927801320278301
$ws.Cells.Item(12, 4).Value = 78.77358394284965
$ws.Cells.Item(12, 5).Value = 300.4003350852541
$ws.Cells.Item(12, 6).Value = 0.4153391009701303

$ws.Cells.Item(13, 1).Value = 'FIRMEZA (N)'
$ws.Cells.Item(13, 2).Value = 'Validação'
$ws.Cells.Item(13, 3).Value = 0.5890851004607165
$ws.Cells.Item(13, 4).Value = 52.4754321071346
$ws.Cells.Item(13, 5).Value = -47.31746203372461
$ws.Cells.Item(13, 6).Value = 1.066119045010663

$ws.Cells.Item(14, 1).Value = 'UBS (%)'
$ws.Cells.Item(14, 2).Value = 'Referência'
$ws.Cells.Item(14, 3).Value = 0.6206853641748344
$ws.Cells.Item(14, 4).Value = 1.937925912849202
$ws.Cells.Item(14, 5).Value = 5.794343485907681
$ws.Cells.Item(14, 6).Value = 0.6206853641748341

$ws.Cells.Item(15, 1).Value = 'UBS (%)'
$ws.Cells.Item(15, 2).Value = 'Predição'
$ws.Cells.Item(15, 3).Value = 0.5837390607208446
$ws.Cells.Item(15, 4).Value = 2.030112941359542
$ws.Cells.Item(15, 5).Value = 6.027106593302635
$ws.Cells.Item(15, 6).Value = 0.605800727325329

$ws.Cells.Item(16, 1).Value = 'UBS (%)'
$ws.Cells.Item(16, 2).Value = 'Validação'
$ws.Cells.Item(16, 3).Value = 0.6862251751601254
$ws.Cells.Item(16, 4).Value = 1.343072576620295
$ws.Cells.Item(16, 5).Value = 1.770590634716598
$ws.Cells.Item(16, 6).Value = 0.8857562505487779
